# Update Wnt10b-Fzd7 NATMI LR-pair TPM results: re-run with new TPM values.
# - Existing "FAPs" sending-cluster rows get refreshed numeric values.
# - A brand-new sending cluster "MuSCs" is added (same 4 target clusters:
#   ECs, FAPs, MuSCs, Resolving-Mac), growing the sheet from 4 to 8 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged, rewritten defensively in case the sheet was blank)
$headers = @(
    "Sending cluster","Ligand symbol","Receptor symbol","Target cluster",
    "Ligand-expressing cells","Ligand detection rate","Ligand average expression value",
    "Ligand total expression value","Ligand derived specificity of average expression value",
    "Ligand derived specificity of total expression value","Receptor-expressing cells",
    "Receptor detection rate","Receptor average expression value","Receptor total expression value",
    "Receptor derived specificity of average expression value","Receptor derived specificity of total expression value",
    "Edge average expression weight","Edge total expression weight",
    "Edge average expression derived specificity","Edge total expression derived specificity"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows: A..D are text (cluster/ligand/receptor labels), E..T are numeric.
$rows = @(
    @("FAPs","Wnt10b","Fzd7","ECs",1,0.3333333333333333,0.08151933333333333,0.244558,0.4155743081351682,0.4155743081351682,3,1,0.8077876666666667,2.423363,0.04902038147436601,0.04902038147436601,0.06585031206155556,0.592652808554,0.02037161111573167,0.02037161111573167),
    @("FAPs","Wnt10b","Fzd7","FAPs",1,0.3333333333333333,0.08151933333333333,0.244558,0.4155743081351682,0.4155743081351682,3,1,7.782607,23.347821,0.4722854529078861,0.4722854529078861,0.6344329342353333,5.709896408117999,0.1962697003344993,0.1962697003344993),
    @("FAPs","Wnt10b","Fzd7","MuSCs",1,0.3333333333333333,0.08151933333333333,0.244558,0.4155743081351682,0.4155743081351682,3,1,5.009781333333333,15.029344,0.3040172587389813,0.3040172587389813,0.408394034439111,3.675546309951999,0.1263417619616025,0.1263417619616025),
    @("FAPs","Wnt10b","Fzd7","Resolving-Mac",1,0.3333333333333333,0.08151933333333333,0.244558,0.4155743081351682,0.4155743081351682,3,1,2.878432333333333,8.635297,0.1746769068787666,0.1746769068787666,0.2346478848584444,2.111830963726,0.07259123472333461,0.07259123472333461),
    @("MuSCs","Wnt10b","Fzd7","ECs",1,0.3333333333333333,0.1146413333333333,0.343924,0.584425691864832,0.584425691864832,3,1,0.8077876666666667,2.423363,0.04902038147436601,0.04902038147436601,0.09260585515688889,0.833452696412,0.02864877035863434,0.02864877035863435),
    @("MuSCs","Wnt10b","Fzd7","FAPs",1,0.3333333333333333,0.1146413333333333,0.343924,0.584425691864832,0.584425691864832,3,1,7.782607,23.347821,0.4722854529078861,0.4722854529078861,0.8922084432893332,8.029875989603999,0.2760157525733868,0.2760157525733868),
    @("MuSCs","Wnt10b","Fzd7","MuSCs",1,0.3333333333333333,0.1146413333333333,0.343924,0.584425691864832,0.584425691864832,3,1,5.009781333333333,15.029344,0.3040172587389813,0.3040172587389813,0.5743280117617777,5.168952105856,0.1776754967773788,0.1776754967773788),
    @("MuSCs","Wnt10b","Fzd7","Resolving-Mac",1,0.3333333333333333,0.1146413333333333,0.343924,0.584425691864832,0.584425691864832,3,1,2.878432333333333,8.635297,0.1746769068787666,0.1746769068787666,0.3299873206031111,2.969885885428,0.102085672155432,0.102085672155432)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
